$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D cells hold text in the source data -- including numeric-looking
# strings such as '598.96' -- which must stay literal text rather than be
# auto-converted to numbers by Excel's smart entry. Force each target cell
# to Text format (@) immediately before writing its new value so only the
# cells actually being rewritten are touched.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.433.62'
$ws.Range('E2').Value = '  -2.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.641.68'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.13'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.80'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.544'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.641.69'
$ws.Range('E9').Value = '  -4.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.146'
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.21'
$ws.Range('E13').Value = '  -3.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.89'
$ws.Range('E14').Value = '  -4.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.120.36'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000185'
$ws.Range('E16').Value = '  -3.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.327.62'
$ws.Range('E17').Value = '  -2.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.672.14'
$ws.Range('E18').Value = '  -3.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.73'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.89'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '364.89'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.40'
$ws.Range('E22').Value = '  -3.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.81'
$ws.Range('E23').Value = '  -4.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  -6.53%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.51'
$ws.Range('E25').Value = '  +5.51%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.82'
$ws.Range('E27').Value = '  -4.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.771.46'
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000103'
$ws.Range('E29').Value = '  -5.12%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '557.86'
$ws.Range('E31').Value = '  -7.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.03'
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.40'
$ws.Range('E33').Value = '  -5.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.93'
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.131'
$ws.Range('E35').Value = '  -3.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').Value = '  -5.97%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.23'
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.43'
$ws.Range('E39').Value = '  -3.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.372'
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('E41').Value = '  -6.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.27'
$ws.Range('E42').Value = '  -5.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.94'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  -8.18%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.16'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₆0303'
$ws.Range('E47').Value = '  -5.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.592'
$ws.Range('E48').Value = '  -3.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '153.75'
$ws.Range('E49').Value = '  -3.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.88'
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.73'
$ws.Range('E51').Value = '  -5.45%  '
